$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-TextOnly($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

Set-TextOnly "D2" "60.166.28"
Set-TextOnly "E2" "  -2.63%  "
Set-TextOnly "D3" "2.386.29"
Set-TextOnly "E3" "  -2.50%  "
Set-TextOnly "E4" "  +0.23%  "
Set-TextValue "D5" "560.46"
Set-TextOnly "E5" "  -2.96%  "
Set-TextValue "D6" "138.25"
Set-TextOnly "E6" "  -2.11%  "
Set-TextOnly "E7" "  -0.08%  "
Set-TextOnly "E8" "  +0.68%  "
Set-TextOnly "D9" "2.387.65"
Set-TextOnly "E9" "  -2.10%  "
Set-TextOnly "E10" "  -4.23%  "
Set-TextOnly "E11" "  -0.98%  "
Set-TextOnly "E12" "  -2.07%  "
Set-TextOnly "E13" "  -1.46%  "
Set-TextOnly "E14" "  -1.81%  "
Set-TextOnly "E15" "  -2.54%  "
Set-TextOnly "E16" "  -3.50%  "
Set-TextOnly "D17" "60.359.03"
Set-TextOnly "E17" "  -2.26%  "
Set-TextValue "D18" "8.42"
Set-TextOnly "D19" "2.385.24"
Set-TextOnly "E19" "  -2.59%  "
Set-TextValue "D20" "10.57"
Set-TextOnly "E20" "  -0.44%  "
Set-TextValue "D21" "324.17"
Set-TextOnly "E21" "  -0.39%  "
Set-TextOnly "E22" "  -1.04%  "
Set-TextValue "D23" "5.97"
Set-TextOnly "E23" "  +0.41%  "
Set-TextOnly "E24" "  -0.13%  "
Set-TextOnly "E25" "  -7.90%  "
Set-TextValue "D26" "64.40"
Set-TextOnly "E26" "  -0.98%  "
Set-TextValue "D27" "552.24"
Set-TextOnly "E27" "  -5.32%  "
Set-TextValue "D28" "7.94"
Set-TextOnly "E28" "  -12.98%  "
Set-TextOnly "D29" "2.506.04"
Set-TextOnly "E29" "  -2.33%  "
Set-TextOnly "D30" "0.0₃0902"
Set-TextOnly "E30" "  -2.42%  "
Set-TextValue "D31" "7.92"
Set-TextOnly "E31" "  -0.21%  "
Set-TextOnly "E32" "  -5.33%  "
Set-TextOnly "E33" "  -3.99%  "
Set-TextOnly "E34" "  -2.51%  "
Set-TextValue "D35" "0.998"
Set-TextOnly "E35" "  -0.46%  "
Set-TextOnly "B36" "Monero"
Set-TextOnly "C36" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D36" "153.26"
Set-TextOnly "E36" "  +1.60%  "
Set-TextOnly "B37" "ImmutableX"
Set-TextOnly "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.41"
Set-TextOnly "E37" "  +0.86%  "
Set-TextValue "D38" "0.367"
Set-TextOnly "E38" "  -1.60%  "
Set-TextOnly "E39" "  -4.38%  "
Set-TextOnly "E40" "  -0.65%  "
Set-TextValue "D41" "5.00"
Set-TextOnly "E41" "  -2.88%  "
Set-TextValue "D43" "41.08"
Set-TextOnly "E43" "  -1.52%  "
Set-TextOnly "E44" "  -2.93%  "
Set-TextOnly "E45" "  -3.98%  "
Set-TextOnly "E46" "  -5.50%  "
Set-TextValue "D47" "142.01"
Set-TextOnly "E47" "  -0.83%  "
Set-TextOnly "E48" "  -2.27%  "
Set-TextOnly "E49" "  -1.89%  "
Set-TextValue "D50" "0.0497"
Set-TextValue "D51" "18.82"
Set-TextOnly "E51" "  -4.16%  "
